# Updated cryptos list with latest Price / Volume(1h) figures (and the
# correct row ordering for LidoStakedEther/USDC and WrappedEther/ShibaInu).
#
# Note: several "Price" strings (column D) look like plain decimal numbers
# (e.g. "570.93", "1.00", "0.0000170"). Excel would normally auto-convert
# such text to a real number on assignment (losing the exact formatting),
# so those are written with a leading apostrophe to force them to stay as
# text, matching the original inline-string cell content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.224.53'
$ws.Range("E2").Value = '  -2.03%  '
$ws.Range("D3").Value = '3.388.56'
$ws.Range("E3").Value = '  -1.54%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'570.93"
$ws.Range("E5").Value = '  -1.47%  '
$ws.Range("D6").Value = "'141.01"
$ws.Range("E6").Value = '  -5.10%  '
$ws.Range("B7").Value = 'LidoStakedEther'
$ws.Range("C7").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D7").Value = '3.384.91'
$ws.Range("E7").Value = '  -1.72%  '
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("E9").Value = '  -0.20%  '
$ws.Range("E10").Value = '  -3.95%  '
$ws.Range("E11").Value = '  -0.84%  '
$ws.Range("D12").Value = "'0.392"
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("D13").Value = '3.968.79'
$ws.Range("E13").Value = '  -1.48%  '
$ws.Range("D14").Value = "'28.14"
$ws.Range("E14").Value = '  +0.45%  '
$ws.Range("E15").Value = '  +1.01%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.390.34'
$ws.Range("E16").Value = '  -1.48%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = "'0.0000170"
$ws.Range("E17").Value = '  -2.64%  '
$ws.Range("D18").Value = '60.443.34'
$ws.Range("E18").Value = '  -1.73%  '
$ws.Range("D19").Value = "'6.27"
$ws.Range("E19").Value = '  -0.86%  '
$ws.Range("D20").Value = "'14.06"
$ws.Range("E20").Value = '  -2.03%  '
$ws.Range("E21").Value = '  -3.74%  '
$ws.Range("D22").Value = "'389.09"
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = "'0.561"
$ws.Range("E23").Value = '  -1.64%  '
$ws.Range("D24").Value = "'73.47"
$ws.Range("E24").Value = '  +1.00%  '
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("E26").Value = '  -4.23%  '
$ws.Range("D27").Value = '3.533.41'
$ws.Range("E27").Value = '  -1.57%  '
$ws.Range("E28").Value = '  -0.45%  '
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("D30").Value = "'7.40"
$ws.Range("E30").Value = '  -5.20%  '
$ws.Range("D31").Value = "'8.03"
$ws.Range("E31").Value = '  -2.80%  '
$ws.Range("D32").Value = "'2.14"
$ws.Range("E32").Value = '  -1.42%  '
$ws.Range("E33").Value = '  -7.10%  '
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").Value = "'23.73"
$ws.Range("E35").Value = '  -1.07%  '
$ws.Range("D36").Value = "'6.94"
$ws.Range("E36").Value = '  -1.92%  '
$ws.Range("D37").Value = '3.420.64'
$ws.Range("E37").Value = '  -1.27%  '
$ws.Range("D38").Value = "'168.12"
$ws.Range("E38").Value = '  +1.32%  '
$ws.Range("E39").Value = '  -6.77%  '
$ws.Range("E40").Value = '  -4.43%  '
$ws.Range("D41").Value = "'0.0777"
$ws.Range("E41").Value = '  -1.96%  '
$ws.Range("D42").Value = "'27.29"
$ws.Range("E42").Value = '  +3.78%  '
$ws.Range("D43").Value = "'0.783"
$ws.Range("E43").Value = '  -1.36%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = "'4.45"
$ws.Range("E45").Value = '  -1.11%  '
$ws.Range("E46").Value = '  -1.80%  '
$ws.Range("D47").Value = "'41.24"
$ws.Range("E47").Value = '  -2.48%  '
$ws.Range("D48").Value = '2.527.39'
$ws.Range("E48").Value = '  -3.35%  '
$ws.Range("E49").Value = '  -3.96%  '
$ws.Range("D50").Value = "'23.24"
$ws.Range("E50").Value = '  +0.38%  '
$ws.Range("D51").Value = "'6.83"
$ws.Range("E51").Value = '  -3.37%  '
